$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.145036666666667
$ws.Range("H2").Value = 3.43511
$ws.Range("I2").Value = 0.4953865629219574
$ws.Range("J2").Value = 0.4953865629219574
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 2.416330038025555
$ws.Range("R2").Value = 21.74697034223
$ws.Range("S2").Value = 0.1898341711725072
$ws.Range("T2").Value = 0.1898341711725072

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.145036666666667
$ws.Range("H3").Value = 3.43511
$ws.Range("I3").Value = 0.4953865629219574
$ws.Range("J3").Value = 0.4953865629219574
$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("Q3").Value = 0.5966251719555555
$ws.Range("R3").Value = 5.369626547599999
$ws.Range("S3").Value = 0.04687267187696965
$ws.Range("T3").Value = 0.04687267187696965

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.145036666666667
$ws.Range("H4").Value = 3.43511
$ws.Range("I4").Value = 0.4953865629219574
$ws.Range("J4").Value = 0.4953865629219574
$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 2.364945372572222
$ws.Range("R4").Value = 21.28450835315
$ws.Range("S4").Value = 0.1857972369690651
$ws.Range("T4").Value = 0.1857972369690651

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.145036666666667
$ws.Range("H5").Value = 3.43511
$ws.Range("I5").Value = 0.4953865629219574
$ws.Range("J5").Value = 0.4953865629219574
$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 0.9276945852144445
$ws.Range("R5").Value = 8.349251266929999
$ws.Range("S5").Value = 0.07288248290341552
$ws.Range("T5").Value = 0.07288248290341552

$ws.Range("G6").Value = 0.4713496666666666
$ws.Range("I6").Value = 0.2039238551060172
$ws.Range("J6").Value = 0.2039238551060172
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.110264333333333
$ws.Range("N6").Value = 6.330793
$ws.Range("O6").Value = 0.3832041185227171
$ws.Range("P6").Value = 0.3832041185227171
$ws.Range("Q6").Value = 0.9946723900952221
$ws.Range("R6").Value = 8.952051510856998
$ws.Range("S6").Value = 0.07814446114165562
$ws.Range("T6").Value = 0.07814446114165562

$ws.Range("G7").Value = 0.4713496666666666
$ws.Range("I7").Value = 0.2039238551060172
$ws.Range("J7").Value = 0.2039238551060172
$ws.Range("O7").Value = 0.0946183755984393
$ws.Range("P7").Value = 0.0946183755984393
$ws.Range("Q7").Value = 0.2455983149822221
$ws.Range("S7").Value = 0.01929494391590285
$ws.Range("T7").Value = 0.01929494391590285

$ws.Range("G8").Value = 0.4713496666666666
$ws.Range("I8").Value = 0.2039238551060172
$ws.Range("J8").Value = 0.2039238551060172
$ws.Range("M8").Value = 2.065388333333333
$ws.Range("N8").Value = 6.196165
$ws.Range("O8").Value = 0.3750550597762889
$ws.Range("P8").Value = 0.3750550597762889
$ws.Range("Q8").Value = 0.9735201024538888
$ws.Range("R8").Value = 8.761680922084999
$ws.Range("S8").Value = 0.07648267366659857
$ws.Range("T8").Value = 0.07648267366659857

$ws.Range("G9").Value = 0.4713496666666666
$ws.Range("I9").Value = 0.2039238551060172
$ws.Range("J9").Value = 0.2039238551060172
$ws.Range("M9").Value = 0.8101876666666666
$ws.Range("N9").Value = 2.430563
$ws.Range("O9").Value = 0.1471224461025547
$ws.Range("P9").Value = 0.1471224461025547
$ws.Range("Q9").Value = 0.3818816866207777
$ws.Range("R9").Value = 3.436935179587
$ws.Range("S9").Value = 0.0300017763818602
$ws.Range("T9").Value = 0.0300017763818602

$ws.Range("G10").Value = 0.6323219999999999
$ws.Range("H10").Value = 1.896966
$ws.Range("I10").Value = 0.2735666300991275
$ws.Range("J10").Value = 0.2735666300991275
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.110264333333333
$ws.Range("N10").Value = 6.330793
$ws.Range("O10").Value = 0.3832041185227171
$ws.Range("P10").Value = 0.3832041185227171
$ws.Range("Q10").Value = 1.334366563782
$ws.Range("R10").Value = 12.009299074038
$ws.Range("S10").Value = 0.1048318593443663
$ws.Range("T10").Value = 0.1048318593443664

$ws.Range("G11").Value = 0.6323219999999999
$ws.Range("H11").Value = 1.896966
$ws.Range("I11").Value = 0.2735666300991275
$ws.Range("J11").Value = 0.2735666300991275
$ws.Range("O11").Value = 0.0946183755984393
$ws.Range("P11").Value = 0.0946183755984393
$ws.Range("Q11").Value = 0.3294734858399999
$ws.Range("R11").Value = 2.96526137256
$ws.Range("S11").Value = 0.02588443015791855
$ws.Range("T11").Value = 0.02588443015791856

$ws.Range("G12").Value = 0.6323219999999999
$ws.Range("H12").Value = 1.896966
$ws.Range("I12").Value = 0.2735666300991275
$ws.Range("J12").Value = 0.2735666300991275
$ws.Range("M12").Value = 2.065388333333333
$ws.Range("N12").Value = 6.196165
$ws.Range("O12").Value = 0.3750550597762889
$ws.Range("P12").Value = 0.3750550597762889
$ws.Range("Q12").Value = 1.30599048171
$ws.Range("R12").Value = 11.75391433539
$ws.Range("S12").Value = 0.1026025488046262
$ws.Range("T12").Value = 0.1026025488046262

$ws.Range("G13").Value = 0.6323219999999999
$ws.Range("H13").Value = 1.896966
$ws.Range("I13").Value = 0.2735666300991275
$ws.Range("J13").Value = 0.2735666300991275
$ws.Range("M13").Value = 0.8101876666666666
$ws.Range("N13").Value = 2.430563
$ws.Range("O13").Value = 0.1471224461025547
$ws.Range("P13").Value = 0.1471224461025547
$ws.Range("Q13").Value = 0.5122994857619999
$ws.Range("R13").Value = 4.610695371857999
$ws.Range("S13").Value = 0.04024779179221641
$ws.Range("T13").Value = 0.04024779179221642

$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.062692
$ws.Range("H14").Value = 0.188076
$ws.Range("I14").Value = 0.02712295187289783
$ws.Range("J14").Value = 0.02712295187289783
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.110264333333333
$ws.Range("N14").Value = 6.330793
$ws.Range("O14").Value = 0.3832041185227171
$ws.Range("P14").Value = 0.3832041185227171
$ws.Range("Q14").Value = 0.1322966915853333
$ws.Range("R14").Value = 1.190670224268
$ws.Range("S14").Value = 0.01039362686418789
$ws.Range("T14").Value = 0.01039362686418789

$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.062692
$ws.Range("H15").Value = 0.188076
$ws.Range("I15").Value = 0.02712295187289783
$ws.Range("J15").Value = 0.02712295187289783
$ws.Range("O15").Value = 0.0946183755984393
$ws.Range("P15").Value = 0.0946183755984393
$ws.Range("Q15").Value = 0.03266587557333333
$ws.Range("R15").Value = 0.29399288016
$ws.Range("S15").Value = 0.00256632964764824
$ws.Range("T15").Value = 0.00256632964764824

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.062692
$ws.Range("H16").Value = 0.188076
$ws.Range("I16").Value = 0.02712295187289783
$ws.Range("J16").Value = 0.02712295187289783
$ws.Range("M16").Value = 2.065388333333333
$ws.Range("N16").Value = 6.196165
$ws.Range("O16").Value = 0.3750550597762889
$ws.Range("P16").Value = 0.3750550597762889
$ws.Range("Q16").Value = 0.1294833253933333
$ws.Range("R16").Value = 1.16534992854
$ws.Range("S16").Value = 0.0101726003359991
$ws.Range("T16").Value = 0.0101726003359991

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.062692
$ws.Range("H17").Value = 0.188076
$ws.Range("I17").Value = 0.02712295187289783
$ws.Range("J17").Value = 0.02712295187289783
$ws.Range("M17").Value = 0.8101876666666666
$ws.Range("N17").Value = 2.430563
$ws.Range("O17").Value = 0.1471224461025547
$ws.Range("P17").Value = 0.1471224461025547
$ws.Range("Q17").Value = 0.05079228519866667
$ws.Range("R17").Value = 0.457130566788
$ws.Range("S17").Value = 0.003990395025062597
$ws.Range("T17").Value = 0.003990395025062597
